$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Concepts")

# Update Date property
$ws1.Range("B8").Value = "2024-09-13T20:57:31+00:00"

# Update Count property
$ws1.Range("B22").NumberFormat = "@"
$ws1.Range("B22").Value = "2"
$ws1.Range("B21").Copy()
$ws1.Range("B22").PasteSpecial(-4122, -4122)

# Copy the formatting of row 2 down to row 3 so new row matches existing style
$ws2.Range("A2:D2").Copy()
$ws2.Range("A3:D3").PasteSpecial(-4122)

# Set the new row values (as text)
$ws2.Range("A3").NumberFormat = "@"
$ws2.Range("A3").Value = "1"
$ws2.Range("A2").Copy()
$ws2.Range("A3").PasteSpecial(-4122, -4122)
$ws2.Range("B3").Value = "unknown"
$ws2.Range("C3").Value = "Unknown"
$ws2.Range("D3").Value = ""
